$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.30423342464539616
$ws.Cells.Item(1, 2).Value = 0.30332923514358612
$ws.Cells.Item(2, 1).Value = -0.25320571575220097
$ws.Cells.Item(2, 2).Value = 0.25005297087934686
$ws.Cells.Item(3, 1).Value = -0.13571060394125389
$ws.Cells.Item(3, 2).Value = 0.13507236979373793
$ws.Cells.Item(4, 1).Value = -0.12707236980188163
$ws.Cells.Item(4, 2).Value = 0.12648895088292811
$ws.Cells.Item(5, 1).Value = -0.12348895088792045
$ws.Cells.Item(5, 2).Value = 0.12149334709188508
$ws.Cells.Item(6, 1).Value = -0.022159077059720644
$ws.Cells.Item(6, 2).Value = 0.021976860748509353
$ws.Cells.Item(7, 1).Value = -0.011976860759947083
$ws.Cells.Item(7, 2).Value = 0.011946515617661468
$ws.Cells.Item(8, 1).Value = -0.016400061700794666
$ws.Cells.Item(8, 2).Value = 0.016306733296272302
$ws.Cells.Item(9, 1).Value = -0.014306733302095864
$ws.Cells.Item(9, 2).Value = 0.014238985674839721
$ws.Cells.Item(10, 1).Value = -0.012238985680809833
$ws.Cells.Item(10, 2).Value = 0.012235099840967933
$ws.Cells.Item(11, 1).Value = -0.0092350998476735668
$ws.Cells.Item(11, 2).Value = 0.0092299653273251181
$ws.Cells.Item(12, 1).Value = -0.0057299653344426993
$ws.Cells.Item(12, 2).Value = 0.0057046021417104242
$ws.Cells.Item(13, 1).Value = -0.0022046021489074974
$ws.Cells.Item(13, 2).Value = 0.0021986425335374093
$ws.Cells.Item(14, 1).Value = 0.005801357455965217
$ws.Cells.Item(14, 2).Value = -0.0058025835349173605
$ws.Cells.Item(15, 1).Value = -0.0080552378990876861
$ws.Cells.Item(15, 2).Value = 0.0080357923176519463
$ws.Cells.Item(16, 1).Value = -0.0060357923237299715
$ws.Cells.Item(16, 2).Value = 0.006003904690976114
$ws.Cells.Item(17, 1).Value = -0.0040039046971411807
$ws.Cells.Item(17, 2).Value = 0.0039999999923390206
$ws.Cells.Item(18, 1).Value = -0.019742078186219203
$ws.Cells.Item(18, 2).Value = 0.019722115374175786
$ws.Cells.Item(19, 1).Value = -0.015722115377438062
$ws.Cells.Item(19, 2).Value = 0.015614053792476934
$ws.Cells.Item(20, 1).Value = -0.0080174877382805931
$ws.Cells.Item(20, 2).Value = 0.0080057787860035035
$ws.Cells.Item(21, 1).Value = -0.004005778789533565
$ws.Cells.Item(21, 2).Value = 0.0039999999964370758
$ws.Cells.Item(22, 1).Value = -0.11034236694274835
$ws.Cells.Item(22, 2).Value = 0.10953734282182293
$ws.Cells.Item(23, 1).Value = -0.040512414581276524
$ws.Cells.Item(23, 2).Value = 0.040101818227140207
$ws.Cells.Item(24, 1).Value = -0.020101818244566694
$ws.Cells.Item(24, 2).Value = 0.019999999982356798
$ws.Cells.Item(25, 1).Value = -0.097334270037721993
$ws.Cells.Item(25, 2).Value = 0.097203639636388672
$ws.Cells.Item(26, 1).Value = -0.094703639642489179
$ws.Cells.Item(26, 2).Value = 0.094536359998574682
$ws.Cells.Item(27, 1).Value = -0.092036360005019358
$ws.Cells.Item(27, 2).Value = 0.09105196778597735
$ws.Cells.Item(28, 1).Value = -0.08905196779330371
$ws.Cells.Item(28, 2).Value = 0.088377789921389116
$ws.Cells.Item(29, 1).Value = -0.081377789933386957
$ws.Cells.Item(29, 2).Value = 0.081181838859471611
$ws.Cells.Item(30, 1).Value = -0.021181838911005446
$ws.Cells.Item(30, 2).Value = 0.021026267350808592
$ws.Cells.Item(31, 1).Value = -0.014026267363652423
$ws.Cells.Item(31, 2).Value = 0.014001891111302456
$ws.Cells.Item(32, 1).Value = -0.0040018911264194656
$ws.Cells.Item(32, 2).Value = 0.0039999999892899041
